$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# Update the date heading
$d.Content.Find.Execute("2025-06-17 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-06-18 Wednesday", 2)

$tbl.Cell(1,1).Range.Text = "9-4="
$tbl.Cell(1,2).Range.Text = "21+43="
$tbl.Cell(1,3).Range.Text = "27+8="
$tbl.Cell(1,4).Range.Text = "82-60="
$tbl.Cell(1,5).Range.Text = "3+75="

$tbl.Cell(2,1).Range.Text = "80-72="
$tbl.Cell(2,2).Range.Text = "29-12="
$tbl.Cell(2,3).Range.Text = "87-15="
$tbl.Cell(2,4).Range.Text = "45+27="
$tbl.Cell(2,5).Range.Text = "59-57="

$tbl.Cell(3,1).Range.Text = "67-57="
$tbl.Cell(3,2).Range.Text = "69-13="
$tbl.Cell(3,3).Range.Text = "17+46="
$tbl.Cell(3,4).Range.Text = "61-60="
$tbl.Cell(3,5).Range.Text = "25+40="

$tbl.Cell(4,1).Range.Text = "83-33="
$tbl.Cell(4,2).Range.Text = "86+13="
$tbl.Cell(4,3).Range.Text = "99-64="
$tbl.Cell(4,4).Range.Text = "83+5="
$tbl.Cell(4,5).Range.Text = "96-12="

$tbl.Cell(5,1).Range.Text = "89-0="
$tbl.Cell(5,2).Range.Text = "58-28="
$tbl.Cell(5,3).Range.Text = "27-24="
$tbl.Cell(5,4).Range.Text = "13+1="
$tbl.Cell(5,5).Range.Text = "75+13="

$tbl.Cell(6,1).Range.Text = "4+16="
$tbl.Cell(6,2).Range.Text = "24+69="
$tbl.Cell(6,3).Range.Text = "46-12="
$tbl.Cell(6,4).Range.Text = "76+8="
$tbl.Cell(6,5).Range.Text = "90-41="

$tbl.Cell(7,1).Range.Text = "32-4="
$tbl.Cell(7,2).Range.Text = "30+48="
$tbl.Cell(7,3).Range.Text = "25+61="
$tbl.Cell(7,4).Range.Text = "74-15="
$tbl.Cell(7,5).Range.Text = "2+66="

$tbl.Cell(8,1).Range.Text = "51+7="
$tbl.Cell(8,2).Range.Text = "92-23="
$tbl.Cell(8,3).Range.Text = "74-73="
$tbl.Cell(8,4).Range.Text = "48-48="
$tbl.Cell(8,5).Range.Text = "10+21="

$tbl.Cell(9,1).Range.Text = "70+12="
$tbl.Cell(9,2).Range.Text = "15+17="
$tbl.Cell(9,3).Range.Text = "81-62="
$tbl.Cell(9,4).Range.Text = "9+72="
$tbl.Cell(9,5).Range.Text = "27-0="

$tbl.Cell(10,1).Range.Text = "2+28="
$tbl.Cell(10,2).Range.Text = "27+11="
$tbl.Cell(10,3).Range.Text = "2+89="
$tbl.Cell(10,4).Range.Text = "2+8="
$tbl.Cell(10,5).Range.Text = "5+56="

$tbl.Cell(11,1).Range.Text = "1+8="
$tbl.Cell(11,2).Range.Text = "71-55="
$tbl.Cell(11,3).Range.Text = "70-53="
$tbl.Cell(11,4).Range.Text = "99-53="
$tbl.Cell(11,5).Range.Text = "16+20="

$tbl.Cell(12,1).Range.Text = "42+29="
$tbl.Cell(12,2).Range.Text = "87-58="
$tbl.Cell(12,3).Range.Text = "0+78="
$tbl.Cell(12,4).Range.Text = "75-46="
$tbl.Cell(12,5).Range.Text = "56-0="

$tbl.Cell(13,1).Range.Text = "41-4="
$tbl.Cell(13,2).Range.Text = "77-23="
$tbl.Cell(13,3).Range.Text = "77-54="
$tbl.Cell(13,4).Range.Text = "15+9="
$tbl.Cell(13,5).Range.Text = "48+36="

$tbl.Cell(14,1).Range.Text = "6+34="
$tbl.Cell(14,2).Range.Text = "75-28="
$tbl.Cell(14,3).Range.Text = "19+78="
$tbl.Cell(14,4).Range.Text = "31-12="
$tbl.Cell(14,5).Range.Text = "86-61="

$tbl.Cell(15,1).Range.Text = "44+37="
$tbl.Cell(15,2).Range.Text = "54+10="
$tbl.Cell(15,3).Range.Text = "73+15="
$tbl.Cell(15,4).Range.Text = "30+60="
$tbl.Cell(15,5).Range.Text = "17+32="

$tbl.Cell(16,1).Range.Text = "12+2="
$tbl.Cell(16,2).Range.Text = "85-76="
$tbl.Cell(16,3).Range.Text = "7+91="
$tbl.Cell(16,4).Range.Text = "80-4="
$tbl.Cell(16,5).Range.Text = "42+26="

$tbl.Cell(17,1).Range.Text = "82+11="
$tbl.Cell(17,2).Range.Text = "92-19="
$tbl.Cell(17,3).Range.Text = "58-34="
$tbl.Cell(17,4).Range.Text = "87-65="
$tbl.Cell(17,5).Range.Text = "29-20="

$tbl.Cell(18,1).Range.Text = "98-34="
$tbl.Cell(18,2).Range.Text = "75-27="
$tbl.Cell(18,3).Range.Text = "90-77="
$tbl.Cell(18,4).Range.Text = "56-24="
$tbl.Cell(18,5).Range.Text = "67+5="

$tbl.Cell(19,1).Range.Text = "73-21="
$tbl.Cell(19,2).Range.Text = "43+26="
$tbl.Cell(19,3).Range.Text = "83-68="
$tbl.Cell(19,4).Range.Text = "93-80="
$tbl.Cell(19,5).Range.Text = "5+14="

$tbl.Cell(20,1).Range.Text = "9-4="
$tbl.Cell(20,2).Range.Text = "50-13="
$tbl.Cell(20,3).Range.Text = "35+6="
$tbl.Cell(20,4).Range.Text = "67-21="
$tbl.Cell(20,5).Range.Text = "39+18="
